# Give slide 1 its own solid background fill (theme color "accent4"),
# instead of inheriting the background from the slide master/layout.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$s.FollowMasterBackground = $false
$s.Background.Fill.Solid()
$s.Background.Fill.ForeColor.ObjectThemeColor = 8
